$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header: "Cap. Percent (%)" in G1
$ws.Range("G1").Value = "Cap. Percent (%)"

# Total of the "abundance*xn" column (D) in D5
$ws.Range("D5").Formula = '=SUM(D2:D4)'

# Per-row percent-of-total-capture column G2:G4
$ws.Range("G2").Formula = '=(D2/$D$5)*100'
$ws.Range("G3:G4").Formula = '=(D3/$D$5)*100'

# Match the column width Excel auto-picked for the new column
$ws.Columns.Item(7).ColumnWidth = 21.6640625

# Move/collapse the selection like the author left it
$ws.Range("G8").Select() | Out-Null
